# Applies the "Automatic update of files" change:
#  - Column C ("Förändrad") date is bumped from 45692 to 45693 for existing data rows (2-36)
#  - Row 36 gains an explicit row height (15, custom)
#  - A new data row 37 is appended ("A 2598-2025")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for all existing data rows.
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 3).Value = 45693
}

# Row 36 picks up an explicit custom row height of 15 in the target file.
$ws.Rows.Item(36).RowHeight = 15

# Append the new row (row 37) with its data.
$ws.Cells.Item(37, 1).Value = "A 2598-2025"

$ws.Cells.Item(37, 2).Value = 45674
$ws.Cells.Item(37, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(37, 3).Value = 45693
$ws.Cells.Item(37, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(37, 4).Value = "OKÄNT"
$ws.Cells.Item(37, 5).Value = "OKÄNT"

$ws.Cells.Item(37, 7).Value = 0.8
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = 0
$ws.Cells.Item(37, 14).Value = 0
$ws.Cells.Item(37, 15).Value = 0
$ws.Cells.Item(37, 16).Value = 0
$ws.Cells.Item(37, 17).Value = 0

# Column R keeps the same (wrap-text) style as the rows above it, but stays blank.
$ws.Range("R36").Copy()
$ws.Range("R37").PasteSpecial(-4122)
